$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cells C2:C10: change the "Förändrad" date serial from 45183 to 45184
foreach ($row in 2..10) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45183) {
        $cell.Value = 45184
    }
}
